$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap on slide 6 ("SOURCES OF FINANCE" table):
#    {FD80E784-CBC0-48E0-A3E2-AE0F3100CCE0} -> {CE6B5B0A-5EDD-40FB-A51D-BA6019FEE7AE}
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{CE6B5B0A-5EDD-40FB-A51D-BA6019FEE7AE}")
    }
}

# ---------------------------------------------------------------------
# 2) Theme colour swap: the deck's two theme parts (the "Integral" theme
#    used by the slide master and the default "Office Theme" used by the
#    notes master) had their colour scheme values swapped. Re-point the
#    slide master's theme colours at the Office Theme palette via the
#    12-slot theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
#    folHlink) - RGB values are stored little-endian (BGR) by COM.
# ---------------------------------------------------------------------
$officeThemeBgr = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeThemeBgr[$i - 1]
}
